$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update the summary figures -------------------------------------------
# VALOR MORA total (E11): 203443 -> 238883
$ws.Range("E11").Value = 238883
# Cant. Trabajadores (C13): 1 -> 2
$ws.Range("C13").Value = 2
# Cant. Periodos (F13): 2 -> 3
$ws.Range("F13").Value = 3

# --- Make room for a new worker row ----------------------------------------
# Shift rows 18 and below down by one row (bounded to the table's columns so
# we don't touch the whole 16384-wide row) so a new data row becomes row 18,
# right below the two existing detail rows and above the blank gap before the
# signature block.
$ws.Range("B18:J18").Insert(-4121) | Out-Null

# The row that used to be the last data row (old row 17, with the "bottom of
# table" border) is now duplicated onto the new row 18; restore row 17 back to
# a normal "middle of table" row by copying row 16's formatting onto it.
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B17:J17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Re-order/re-populate the worker detail rows ---------------------------
# Row 16: CC / 1047459833 / IVETH CAROLINA MARRUGO PAUTT / period 2408
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047459833"
$ws.Range("D16").Value = "IVETH CAROLINA MARRUGO PAUTT"
$ws.Range("E16").Value = "2408"
$ws.Range("F16").Value = 132680
$ws.Range("G16").Value = 3317000

# Row 17: CC / 1047459833 / IVETH CAROLINA MARRUGO PAUTT / period 2407
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047459833"
$ws.Range("D17").Value = "IVETH CAROLINA MARRUGO PAUTT"
$ws.Range("E17").Value = "2407"
$ws.Range("F17").Value = 70763
$ws.Range("G17").Value = 3317000

# Row 18 (new worker): CC / 1143396862 / NEY SMITH CERVANTES BOLAÑOS / period 2006
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143396862"
$ws.Range("D18").Value = "NEY SMITH CERVANTES BOLAÑOS"
$ws.Range("E18").Value = "2006"
$ws.Range("F18").Value = 35440
$ws.Range("G18").Value = 886000

# --- Re-fit the bestFit columns ---------------------------------------------
$ws.Columns("B:J").AutoFit() | Out-Null
